$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in column C grades for each student
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 2.5
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0

# Update selection to reflect the new active cell (C7)
$ws.Range("C7").Select()
